# UTS - Perbaikan error/bug
# Template barang: replace the seeded sample rows with a single corrected
# entry and bump harga_beli/harga_jual for row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fix product code/name + update the prices.
$ws.Range("B2").Value = "SBW"
$ws.Range("C2").Value = "Sabun Wajah"
$ws.Range("D2").Value = 22300
$ws.Range("E2").Value = 25300

# Rows 3-6 no longer hold sample data - wipe both content and formatting
# so the cells fall back to the workbook's default (unstyled) state.
$clearRange = $ws.Range("A3:E6")
$clearRange.ClearContents()
$clearRange.Style = "Normal"

# Move the active selection to F4, matching where the editor left off.
$ws.Range("F4").Select()
